$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null
$ws.Range("H132").Value = 2621.1365
$ws.Range("I132").Value = 2024.0238
$ws.Range("J132").Value = 15160.5
$ws.Range("K132").Value = 6072.0714
$ws.Range("L132").Value = 45481.5
$ws.Range("M132").Value = -3542.0714
$ws.Range("N132").Value = -50541.5
$ws.Range("H141").Value = 4850.2
$ws.Range("I141").Value = 5381.077
$ws.Range("K141").Value = 16143.231
$ws.Range("M141").Value = -10963.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27816.453
$ws.Range("I32").Value = 29007.658
$ws.Range("J32").Value = 16500
$ws.Range("K32").Value = 29007.658
$ws.Range("L32").Value = 16500
$ws.Range("M32").Value = -28720.658
$ws.Range("N32").Value = -17074
$ws.Range("H45").Value = 2784.6
$ws.Range("I45").Value = 1014.9231
$ws.Range("K45").Value = 1014.9231
$ws.Range("M45").Value = -637.9231
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 419981.5
$ws.Range("J141").Value = 419981.5
$ws.Range("L141").Value = 419981.5
$ws.Range("N141").Value = -430341.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1436440.8
$ws.Range("I4").Value = 1009652
$ws.Range("K4").Value = 3028956
$ws.Range("M4").Value = -3028844
$ws.Range("H17").Value = 1052.25
$ws.Range("I17").Value = 105
$ws.Range("K17").Value = 315
$ws.Range("M17").Value = -146
$ws.Range("H25").Value = 7566.6665
$ws.Range("I25").Value = 700
$ws.Range("K25").Value = 2100
$ws.Range("M25").Value = -1931
$ws.Range("H30").Value = 7566.6665
$ws.Range("I30").Value = 700
$ws.Range("K30").Value = 2100
$ws.Range("M30").Value = -1998
$ws.Range("H36").Value = 772.6667
$ws.Range("I36").Value = 772.6667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2318.0001
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2149.0001
$ws.Range("N36").Value = $null
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = $null
$ws.Range("H48").Value = 4166.3335
$ws.Range("J48").Value = 4599.6
$ws.Range("L48").Value = 13798.8
$ws.Range("N48").Value = -14298.8
$ws.Range("H61").Value = 279
$ws.Range("J61").Value = 416.66666
$ws.Range("L61").Value = 1249.99998
$ws.Range("N61").Value = -1679.99998
$ws.Range("H62").Value = 5603.25
$ws.Range("J62").Value = 5603.25
$ws.Range("L62").Value = 16809.75
$ws.Range("N62").Value = -18181.75
$ws.Range("H63").Value = 400
$ws.Range("I63").Value = 400
$ws.Range("K63").Value = 1200
$ws.Range("M63").Value = -451
$ws.Range("H64").Value = 10999
$ws.Range("J64").Value = 10999
$ws.Range("L64").Value = 32997
$ws.Range("N64").Value = -33537
$ws.Range("H65").Value = 5603.25
$ws.Range("J65").Value = 5603.25
$ws.Range("L65").Value = 50429.25
$ws.Range("N65").Value = -57293.25
$ws.Range("H66").Value = 400
$ws.Range("I66").Value = 400
$ws.Range("K66").Value = 3600
$ws.Range("M66").Value = 144
$ws.Range("H67").Value = 10999
$ws.Range("J67").Value = 10999
$ws.Range("L67").Value = 32997
$ws.Range("N67").Value = -34869
$ws.Range("H70").Value = 20000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null
$ws.Range("H73").Value = 20000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null
$ws.Range("H75").Value = 950
$ws.Range("J75").Value = 950
$ws.Range("L75").Value = 2850
$ws.Range("N75").Value = -4846
$ws.Range("I76").Value = 3013
$ws.Range("J76").Value = 20000
$ws.Range("K76").Value = 9039
$ws.Range("L76").Value = 60000
$ws.Range("M76").Value = -8656
$ws.Range("N76").Value = -60766
$ws.Range("H78").Value = 950
$ws.Range("J78").Value = 950
$ws.Range("L78").Value = 8550
$ws.Range("N78").Value = -18534
$ws.Range("I79").Value = 3013
$ws.Range("J79").Value = 20000
$ws.Range("K79").Value = 9039
$ws.Range("L79").Value = 60000
$ws.Range("M79").Value = -7713
$ws.Range("N79").Value = -62652
$ws.Range("H80").Value = 5887
$ws.Range("J80").Value = 6591.3335
$ws.Range("L80").Value = 19774.0005
$ws.Range("N80").Value = -21646.0005
$ws.Range("H83").Value = 5887
$ws.Range("J83").Value = 6591.3335
$ws.Range("L83").Value = 59322.0015
$ws.Range("N83").Value = -68682.0015
$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("K87").Value = 3000
$ws.Range("M87").Value = -1752
$ws.Range("H88").Value = 11124.625
$ws.Range("J88").Value = 11428.571
$ws.Range("L88").Value = 34285.713
$ws.Range("N88").Value = -35141.713
$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("K90").Value = 9000
$ws.Range("M90").Value = -2760
$ws.Range("H91").Value = 11124.625
$ws.Range("J91").Value = 11428.571
$ws.Range("L91").Value = 34285.713
$ws.Range("N91").Value = -37249.713
$ws.Range("H97").Value = 406.58334
$ws.Range("I97").Value = 228
$ws.Range("J97").Value = 466.1111
$ws.Range("K97").Value = 684
$ws.Range("L97").Value = 1398.3333
$ws.Range("M97").Value = -188
$ws.Range("N97").Value = -2390.3333
$ws.Range("H98").Value = 4492.4
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 4492.4
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 13477.2
$ws.Range("M98").Value = $null
$ws.Range("N98").Value = -16473.2
$ws.Range("H99").Value = 6035.7
$ws.Range("J99").Value = 6638.8335
$ws.Range("L99").Value = 19916.5005
$ws.Range("N99").Value = -24408.5005
$ws.Range("H100").Value = 3000
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 9000
$ws.Range("N100").Value = -10622
$ws.Range("H104").Value = 3066.6667
$ws.Range("I104").Value = 2500
$ws.Range("J104").Value = 3350
$ws.Range("K104").Value = 7500
$ws.Range("L104").Value = 10050
$ws.Range("M104").Value = -4879
$ws.Range("N104").Value = -15292
$ws.Range("H105").Value = 9110.223
$ws.Range("J105").Value = 9110.223
$ws.Range("L105").Value = 27330.669
$ws.Range("N105").Value = -32572.669
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 60000
$ws.Range("N112").Value = -62216

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2389
$ws.Range("I102").Value = 1453.9166
$ws.Range("K102").Value = 1453.9166
$ws.Range("M102").Value = 168.0834
$ws.Range("H107").Value = 44572.957
$ws.Range("I107").Value = 59687.53
$ws.Range("J107").Value = 1748.3334
$ws.Range("K107").Value = 59687.53
$ws.Range("L107").Value = 1748.3334
$ws.Range("M107").Value = -57767.53
$ws.Range("N107").Value = -5588.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3516.25
$ws.Range("I40").Value = 2738.3333
$ws.Range("J40").Value = 5850
$ws.Range("K40").Value = 2738.3333
$ws.Range("L40").Value = 5850
$ws.Range("M40").Value = -2602.3333
$ws.Range("N40").Value = -6122
$ws.Range("H68").Value = 3097.1428
$ws.Range("I68").Value = 1850
$ws.Range("J68").Value = 3596
$ws.Range("K68").Value = 1850
$ws.Range("L68").Value = 3596
$ws.Range("M68").Value = -1101
$ws.Range("N68").Value = -5094
$ws.Range("H71").Value = 3097.1428
$ws.Range("I71").Value = 1850
$ws.Range("J71").Value = 3596
$ws.Range("K71").Value = 9250
$ws.Range("L71").Value = 17980
$ws.Range("M71").Value = -5506
$ws.Range("N71").Value = -25468
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 14950
$ws.Range("J22").Value = 14950
$ws.Range("L22").Value = 14950
$ws.Range("N22").Value = -15536
$ws.Range("H113").Value = 938.48
$ws.Range("I113").Value = 585.6429000000001
$ws.Range("K113").Value = 1756.9287
$ws.Range("M113").Value = 413.0712999999998
$ws.Range("H122").Value = 1510.0625
$ws.Range("I122").Value = 1497.2142
$ws.Range("K122").Value = 4491.642599999999
$ws.Range("M122").Value = -2041.642599999999
$ws.Range("H125").Value = 53107.5
$ws.Range("J125").Value = 53107.5
$ws.Range("L125").Value = 53107.5
$ws.Range("N125").Value = -62947.5
$ws.Range("H136").Value = 2209.8965
$ws.Range("I136").Value = 1842.0834
$ws.Range("J136").Value = 3975.4
$ws.Range("K136").Value = 5526.2502
$ws.Range("L136").Value = 11926.2
$ws.Range("M136").Value = -2976.2502
$ws.Range("N136").Value = -17026.2
